$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = '@'
    $r.Value = $val
    $r.Style = 'Normal'
}

$ws.Range('D2').Value = '25.682.91'
$ws.Range('E2').Value = '  -3.59%  '
$ws.Range('D3').Value = '1.744.93'
$ws.Range('E3').Value = '  -5.62%  '
Set-TextValue 'D4' '1.001'
$ws.Range('E4').Value = '  -0.05%  '
Set-TextValue 'D5' '236.24'
$ws.Range('E5').Value = '  -10.48%  '
$ws.Range('E6').Value = '  -0.07%  '
Set-TextValue 'D7' '0.4926'
$ws.Range('E7').Value = '  -7.68%  '
Set-TextValue 'D8' '41.59'
$ws.Range('E8').Value = '  -7.79%  '
Set-TextValue 'D9' '0.2518'
$ws.Range('E9').Value = '  -20.72%  '
Set-TextValue 'D10' '0.06021'
$ws.Range('E10').Value = '  -13.54%  '
$ws.Range('D11').Value = '1.745.29'
$ws.Range('E11').Value = '  -5.87%  '
Set-TextValue 'D12' '0.06832'
$ws.Range('E12').Value = '  -12.68%  '
Set-TextValue 'D13' '14.84'
$ws.Range('E13').Value = '  -21.50%  '
Set-TextValue 'D14' '4.452'
$ws.Range('E14').Value = '  -11.86%  '
Set-TextValue 'D15' '76.78'
$ws.Range('E15').Value = '  -14.42%  '
Set-TextValue 'D16' '0.5657'
$ws.Range('E16').Value = '  -26.68%  '
$ws.Range('E17').Value = '  -0.06%  '
$ws.Range('E18').Value = '  +0.01%  '
$ws.Range('D19').Value = '25.733.00'
$ws.Range('E19').Value = '  -3.52%  '
Set-TextValue 'D20' '11.27'
$ws.Range('E20').Value = '  -20.32%  '
Set-TextValue 'D21' '0.000006552'
$ws.Range('E21').Value = '  -18.01%  '
$ws.Range('D22').Value = '1.964.75'
$ws.Range('E22').Value = '  -6.01%  '
Set-TextValue 'D23' '4.000'
$ws.Range('E23').Value = '  -14.04%  '
$ws.Range('B24').Value = 'Cosmos'
$ws.Range('C24').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue 'D24' '7.876'
$ws.Range('E24').Value = '  -16.15%  '
$ws.Range('B25').Value = 'Chainlink'
$ws.Range('C25').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue 'D25' '5.011'
$ws.Range('E25').Value = '  -16.96%  '
Set-TextValue 'D26' '136.50'
$ws.Range('E26').Value = '  -4.31%  '
Set-TextValue 'D27' '1.479'
$ws.Range('E27').Value = '  -13.30%  '
Set-TextValue 'D28' '1.815'
$ws.Range('E28').Value = '  -18.06%  '
Set-TextValue 'D29' '14.65'
$ws.Range('E29').Value = '  -14.55%  '
Set-TextValue 'D30' '101.90'
$ws.Range('E30').Value = '  -8.82%  '
Set-TextValue 'D31' '3.757'
$ws.Range('E31').Value = '  -13.23%  '
Set-TextValue 'D32' '0.07967'
$ws.Range('E32').Value = '  -9.14%  '
Set-TextValue 'D33' '3.371'
$ws.Range('E33').Value = '  -18.07%  '
Set-TextValue 'D34' '0.04377'
$ws.Range('E34').Value = '  -10.07%  '
Set-TextValue 'D35' '1.000'
$ws.Range('E35').Value = '  -0.04%  '
Set-TextValue 'D36' '2.631'
$ws.Range('E36').Value = '  -8.88%  '
Set-TextValue 'D37' '0.9709'
$ws.Range('E37').Value = '  -14.83%  '
Set-TextValue 'D38' '0.6034'
$ws.Range('E38').Value = '  -18.53%  '
Set-TextValue 'D39' '2.675'
$ws.Range('E39').Value = '  -13.87%  '
Set-TextValue 'D40' '2.008'
$ws.Range('E40').Value = '  -14.93%  '
$ws.Range('E41').Value = '  -0.05%  '
Set-TextValue 'D42' '102.68'
$ws.Range('E42').Value = '  -6.10%  '
Set-TextValue 'D43' '0.01499'
$ws.Range('E43').Value = '  -14.17%  '
Set-TextValue 'D44' '0.7551'
$ws.Range('E44').Value = '  -16.81%  '
Set-TextValue 'D45' '5.173'
$ws.Range('E45').Value = '  -12.55%  '
Set-TextValue 'D46' '0.3725'
$ws.Range('E46').Value = '  -23.00%  '
Set-TextValue 'D47' '0.05271'
$ws.Range('E47').Value = '  -9.60%  '
Set-TextValue 'D48' '0.1061'
$ws.Range('E48').Value = '  -15.23%  '
Set-TextValue 'D49' '30.00'
$ws.Range('E49').Value = '  -14.81%  '
Set-TextValue 'D50' '5.886'
$ws.Range('E50').Value = '  -23.70%  '
Set-TextValue 'D51' '52.24'
$ws.Range('E51').Value = '  -13.68%  '
